$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric / identifier fields -------------------------------------------------
$ws.Range("A2").Value = 102077473
$ws.Range("B2").Value = 96367
$ws.Range("E2").Value = 219874

# --- Species identification fields ------------------------------------------------
$ws.Range("F2").Value = "Nattviol"
$ws.Range("G2").Value = "Platanthera bifolia"
$ws.Range("H2").Value = "(L.) Rich."

# Antal (I2): becomes text "3" (was an empty text cell before)
$ws.Range("I2").Value = "'3"
$ws.Range("I2").ClearFormats()

# Enhet / Alder-Stadium / Kon / Metod (J2, K2, L2, N2): newly-created, blank text cells
$ws.Range("J2").Value = "'"
$ws.Range("J2").ClearFormats()
$ws.Range("K2").Value = "'"
$ws.Range("K2").ClearFormats()
$ws.Range("L2").Value = "'"
$ws.Range("L2").ClearFormats()
$ws.Range("N2").Value = "'"
$ws.Range("N2").ClearFormats()

# --- Location fields ---------------------------------------------------------------
$ws.Range("P2").Value = "Tvetaspåret, Tveta, Srm"
$ws.Range("Q2").Value = 647720.9098417715
$ws.Range("R2").Value = 6560694.968483768
$ws.Range("S2").Value = 10

# --- Dates (kept as literal text, matching the source workbook's storage) ------------
$ws.Range("Y2").Value = "'2022-06-28"
$ws.Range("Y2").ClearFormats()
$ws.Range("AA2").Value = "'2022-07-05"
$ws.Range("AA2").ClearFormats()

# Bestamningsmetod (AF2): newly-created, blank text cell
$ws.Range("AF2").Value = "'"
$ws.Range("AF2").ClearFormats()

# Biotop-beskrivning (AI2): removed entirely
$ws.Range("AI2").ClearContents()

# --- People ----------------------------------------------------------------------------
$ws.Range("AW2").Value = "Åsa Johansson"
$ws.Range("AX2").Value = "Åsa Johansson"

Write-Output "edit applied"
